$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @(
    728.75,
    611,
    447,
    867.5,
    811.25,
    877.25,
    789.25,
    869.5,
    882.75,
    704.25,
    593.25,
    759.75,
    661,
    720.25,
    510,
    821.5,
    752.5,
    632.25,
    642.5,
    719,
    673.25,
    668,
    766.75,
    851.25,
    759,
    709.75,
    663.5,
    829
)

$row = 2
foreach ($v in $values) {
    $ws.Cells.Item($row, 2).Value = $v
    $row++
}
